$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2D training schedule values (rows 2-6, columns B-H)
$data = @{
    2 = @{ B=4; C=3; D=5; E=8; F=1; G=5; H=56 }
    3 = @{ B=4; C=4; D=6; E=8; F=2; G=4; H=45 }
    4 = @{ B=0; C=1; D=3; E=4; F=3; G=3; H=34 }
    5 = @{ B=4; C=1; D=8; E=3; F=4; G=2; H=23 }
    6 = @{ B=1; C=0; D=6; E=1; F=5; G=1; H=12 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

# Set the active cell selection to I1 (no break screen column)
$ws.Range("I1").Select()
